$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates per the crypto price/volume refresh.
# Numeric-looking price strings must be forced to remain plain text
# (matching the source workbook, where Price/Volume columns are stored
# as literal strings, not numbers) by briefly switching the cell to a
# text format before assigning, then clearing the format override so
# the cell keeps its original (default) style.

$ws.Range("D2").Value = "43.148.18"
$ws.Range("E2").Value = "  +4.54%  "
$ws.Range("D3").Value = "2.252.17"
$ws.Range("E3").Value = "  +3.67%  "
$ws.Range("E4").Value = "  +0.00%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "245.02"
$c.ClearFormats()
$ws.Range("E5").Value = "  +3.34%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.618"
$c.ClearFormats()
$ws.Range("E6").Value = "  +1.70%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "75.99"
$c.ClearFormats()
$ws.Range("E7").Value = "  +9.36%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("E9").Value = "  +7.12%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "41.39"
$c.ClearFormats()
$ws.Range("E10").Value = "  +5.59%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0937"
$c.ClearFormats()
$ws.Range("E11").Value = "  +2.10%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "7.03"
$c.ClearFormats()
$ws.Range("E12").Value = "  +5.15%  "
$ws.Range("E13").Value = "  +1.14%  "
$ws.Range("D14").Value = "2.590.97"
$ws.Range("E14").Value = "  +3.76%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "14.60"
$c.ClearFormats()
$ws.Range("E15").Value = "  +4.86%  "
$ws.Range("D16").Value = "2.256.24"
$ws.Range("E16").Value = "  +4.36%  "
$ws.Range("E17").Value = "  +1.93%  "
$ws.Range("D18").Value = "43.089.69"
$ws.Range("E18").Value = "  +4.92%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.0000106"
$c.ClearFormats()
$ws.Range("E19").Value = "  +5.32%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "71.49"
$c.ClearFormats()
$ws.Range("E20").Value = "  +1.57%  "
$ws.Range("E21").Value = "  +2.70%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "10.04"
$c.ClearFormats()
$ws.Range("E22").Value = "  +7.60%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "230.76"
$c.ClearFormats()
$ws.Range("E23").Value = "  +2.65%  "
$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.20"
$c.ClearFormats()
$ws.Range("E24").Value = "  +16.71%  "
$ws.Range("E25").Value = "  -0.06%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "10.98"
$c.ClearFormats()
$ws.Range("E26").Value = "  +2.74%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "3.52"
$c.ClearFormats()
$ws.Range("E27").Value = "  +1.97%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.31"
$c.ClearFormats()
$ws.Range("E28").Value = "  +6.52%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "39.81"
$c.ClearFormats()
$ws.Range("E29").Value = "  +33.44%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "2.25"
$c.ClearFormats()
$ws.Range("E30").Value = "  +2.88%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "174.26"
$c.ClearFormats()
$ws.Range("E31").Value = "  +4.16%  "
$ws.Range("E32").Value = "  +2.83%  "
$ws.Range("E33").Value = "  +5.35%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.37"
$c.ClearFormats()
$ws.Range("E34").Value = "  +5.62%  "
$ws.Range("E35").Value = "  +2.05%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.112"
$c.ClearFormats()
$ws.Range("E36").Value = "  +11.43%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "4.34"
$c.ClearFormats()
$ws.Range("E37").Value = "  +7.16%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.0334"
$c.ClearFormats()
$ws.Range("E38").Value = "  +19.05%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "13.32"
$c.ClearFormats()
$ws.Range("E39").Value = "  +14.33%  "
$ws.Range("E40").Value = "  +4.59%  "
$ws.Range("E41").Value = "  +3.90%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.205"
$c.ClearFormats()
$ws.Range("E42").Value = "  +8.60%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "60.20"
$c.ClearFormats()
$ws.Range("E43").Value = "  +2.63%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "106.21"
$c.ClearFormats()
$ws.Range("E44").Value = "  +9.68%  "
$ws.Range("E45").Value = "  +6.11%  "
$ws.Range("B46").Value = "WOONetwork"
$ws.Range("C46").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.485"
$c.ClearFormats()
$ws.Range("E46").Value = "  +31.31%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.0998"
$c.ClearFormats()
$ws.Range("E47").Value = "  +2.93%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.42"
$c.ClearFormats()
$ws.Range("E48").Value = "  +10.93%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.12"
$c.ClearFormats()
$ws.Range("E49").Value = "  +4.18%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.16"
$c.ClearFormats()
$ws.Range("E50").Value = "  +3.90%  "
$ws.Range("D51").Value = "2.465.19"
$ws.Range("E51").Value = "  +3.90%  "
